# Actualización automática 2025-06-24 08:30:09
# Updates June ("junio") sales figures and the dependent summary sheets
# for GUERRERO FAREZ FABIAN MAURICIO.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: VENTAS POR GRUPO (per-client sales by product group)
# ---------------------------------------------------------------------------
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

$wsGrupo.Range("D5").Value  = 739.58
$wsGrupo.Range("D14").Value = 570.24
$wsGrupo.Range("H16").Value = 980.1
$wsGrupo.Range("O16").Value = 3248.09
$wsGrupo.Range("P16").Value = 379.57
$wsGrupo.Range("D29").Value = 2809.72
$wsGrupo.Range("L29").Value = 1677.15
$wsGrupo.Range("Q29").Value = 364.73

# Weekly "clients with sales" counters on the totals row
$wsGrupo.Range("D54").Value = "9 de 52"
$wsGrupo.Range("H54").Value = "1 de 52"
$wsGrupo.Range("L54").Value = "5 de 52"
$wsGrupo.Range("P54").Value = "5 de 52"
$wsGrupo.Range("Q54").Value = "2 de 52"

# ---------------------------------------------------------------------------
# Sheet: VENTA MENSUAL (monthly sales - "junio" column F)
# ---------------------------------------------------------------------------
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F5").Value  = 3500.52
$wsMensual.Range("F14").Value = 1571.01
$wsMensual.Range("F16").Value = 5153.05
$wsMensual.Range("F29").Value = 5400.91
$wsMensual.Range("F54").Value = 53454.81

# ---------------------------------------------------------------------------
# Sheet: CUMPLIMIENTO MENSUAL (budget vs. actual compliance)
# ---------------------------------------------------------------------------
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$wsCumplimiento.Range("D3").Value = 10789.46
$wsCumplimiento.Range("E3").Value = 16667.5476
$wsCumplimiento.Range("F3").Value = 0.392958335343142

$wsCumplimiento.Range("D7").Value = 980.1
$wsCumplimiento.Range("E7").Value = 419.9
$wsCumplimiento.Range("F7").Value = 0.7000714285714286

$wsCumplimiento.Range("D10").Value = 653.01
$wsCumplimiento.Range("E10").Value = 647.49
$wsCumplimiento.Range("F10").Value = 0.5021222606689735

$wsCumplimiento.Range("D14").Value = 670.46
$wsCumplimiento.Range("E14").Value = 295.54
$wsCumplimiento.Range("F14").Value = 0.6940579710144927

$wsCumplimiento.Range("D15").Value = 19281.21
$wsCumplimiento.Range("E15").Value = -5781.209999999999
$wsCumplimiento.Range("F15").Value = 1.428237777777778

$wsCumplimiento.Range("D18").Value = 6387.11
$wsCumplimiento.Range("E18").Value = -3187.11
$wsCumplimiento.Range("F18").Value = 1.995971875

$wsCumplimiento.Range("D19").Value = 53454.81
$wsCumplimiento.Range("E19").Value = 40992.63064517915
$wsCumplimiento.Range("F19").Value = 0.5659741506476542
